# Apply the changes described by the diff:
# 1. Update the Date value on the Metadata sheet.
# 2. Swap the "AK" and "AL" columns (header text, data values, and column widths)
#    on the Elements sheet (a new "Mapping: business spec -> ROR DropZone extension"
#    mapping column is inserted before the existing "Mapping: RIM Mapping" column).

$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: update Date ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# --- 2. Elements sheet: swap columns AK (37) and AL (38) ---
$wsEl = $wb.Worksheets.Item("Elements")

# Determine the last used row on the Elements sheet.
$lastRow = $wsEl.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $wsEl.Cells.Item($r, 37)
    $alCell = $wsEl.Cells.Item($r, 38)

    $akVal = $akCell.Value()
    $alVal = $alCell.Value()

    if ($akVal -ne $alVal) {
        $akCell.Value = $alVal
        $alCell.Value = $akVal
    }
}

# Swap the column widths too, since the wider column (the new mapping text)
# moved from AL to AK and vice versa. (Column AK should become as wide as the
# old AL column - about 65.18 characters - and AL should shrink to the old AK
# width of about 24.98 characters.)
$wsEl.Columns.Item(37).ColumnWidth = 64.33333333333333
$wsEl.Columns.Item(38).ColumnWidth = 24.166666666666668
